$wb = $excel.ActiveWorkbook

# ALC sheet (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 2642.6316
$ws.Cells.Item(17, 10).Value = 2313.7222
$ws.Cells.Item(17, 12).Value = 6941.1666
$ws.Cells.Item(17, 14).Value = -7277.1666
$ws.Cells.Item(43, 8).Value = 804.1667
$ws.Cells.Item(43, 9).Value = 815
$ws.Cells.Item(43, 10).Value = 750
$ws.Cells.Item(43, 11).Value = 815
$ws.Cells.Item(43, 12).Value = 750
$ws.Cells.Item(43, 13).Value = -746
$ws.Cells.Item(43, 14).Value = -888
$ws.Cells.Item(53, 8).Value = 36676.332
$ws.Cells.Item(53, 9).Value = 36676.332
$ws.Cells.Item(53, 11).Value = 36676.332
$ws.Cells.Item(53, 13).Value = -36039.332
$ws.Cells.Item(106, 8).Value = 5623.75
$ws.Cells.Item(106, 9).Value = 5623.75
$ws.Cells.Item(106, 11).Value = 5623.75
$ws.Cells.Item(106, 13).Value = -4992.75
$ws.Cells.Item(132, 8).Value = 1350.0416
$ws.Cells.Item(132, 9).Value = 1352.2
$ws.Cells.Item(132, 10).Value = 1346.4445
$ws.Cells.Item(132, 11).Value = 4056.6
$ws.Cells.Item(132, 12).Value = 4039.3335
$ws.Cells.Item(132, 13).Value = -1526.6
$ws.Cells.Item(132, 14).Value = -9099.333500000001
$ws.Cells.Item(137, 8).Value = 2404.3845
$ws.Cells.Item(137, 9).Value = 1705.5555
$ws.Cells.Item(137, 11).Value = 5116.666499999999
$ws.Cells.Item(137, 13).Value = -2566.666499999999
$ws.Cells.Item(138, 8).Value = 2126.4407
$ws.Cells.Item(138, 10).Value = 2525.0435
$ws.Cells.Item(138, 12).Value = 7575.130500000001
$ws.Cells.Item(138, 14).Value = -17855.1305

# ARM sheet (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 4579.72
$ws.Cells.Item(32, 9).Value = 3488.2327
$ws.Cells.Item(32, 10).Value = 11284.571
$ws.Cells.Item(32, 11).Value = 3488.2327
$ws.Cells.Item(32, 12).Value = 11284.571
$ws.Cells.Item(32, 13).Value = -3201.2327
$ws.Cells.Item(32, 14).Value = -11858.571
$ws.Cells.Item(61, 8).Value = 2100.3333
$ws.Cells.Item(61, 9).Value = 2101.4285
$ws.Cells.Item(61, 10).Value = 2098.8
$ws.Cells.Item(61, 11).Value = 2101.4285
$ws.Cells.Item(61, 12).Value = 2098.8
$ws.Cells.Item(61, 13).Value = -1889.4285
$ws.Cells.Item(61, 14).Value = -2522.8
$ws.Cells.Item(74, 8).Value = 1140.8857
$ws.Cells.Item(74, 9).Value = 583.5
$ws.Cells.Item(74, 11).Value = 583.5
$ws.Cells.Item(74, 13).Value = 290.5
$ws.Cells.Item(77, 8).Value = 1140.8857
$ws.Cells.Item(77, 9).Value = 583.5
$ws.Cells.Item(77, 11).Value = 2917.5
$ws.Cells.Item(77, 13).Value = 1450.5
$ws.Cells.Item(122, 8).Value = 1496.6666
$ws.Cells.Item(122, 9).Value = 1026.6666
$ws.Cells.Item(122, 10).Value = 1966.6666
$ws.Cells.Item(122, 11).Value = 3079.9998
$ws.Cells.Item(122, 12).Value = 5899.9998
$ws.Cells.Item(122, 13).Value = -629.9998000000001
$ws.Cells.Item(122, 14).Value = -10799.9998
$ws.Cells.Item(132, 8).Value = 1544.5555
$ws.Cells.Item(132, 9).Value = 1410.5714
$ws.Cells.Item(132, 11).Value = 4231.7142
$ws.Cells.Item(132, 13).Value = -1701.7142
$ws.Cells.Item(136, 8).Value = 2100.3333
$ws.Cells.Item(136, 9).Value = 2101.4285
$ws.Cells.Item(136, 10).Value = 2098.8
$ws.Cells.Item(136, 11).Value = 6304.2855
$ws.Cells.Item(136, 12).Value = 6296.400000000001
$ws.Cells.Item(136, 13).Value = -3754.2855
$ws.Cells.Item(136, 14).Value = -11396.4

# BSM sheet (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 6924.522
$ws.Cells.Item(134, 9).Value = 9505.6
$ws.Cells.Item(134, 10).Value = 2085
$ws.Cells.Item(134, 11).Value = 28516.8
$ws.Cells.Item(134, 12).Value = 6255
$ws.Cells.Item(134, 13).Value = -25981.8
$ws.Cells.Item(134, 14).Value = -11325

# CRP sheet (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2857.5881
$ws.Cells.Item(31, 9).Value = 2841.111
$ws.Cells.Item(31, 10).Value = 2876.125
$ws.Cells.Item(31, 11).Value = 2841.111
$ws.Cells.Item(31, 12).Value = 2876.125
$ws.Cells.Item(31, 13).Value = -2546.111
$ws.Cells.Item(31, 14).Value = -3466.125
$ws.Cells.Item(34, 8).Value = 2857.5881
$ws.Cells.Item(34, 9).Value = 2841.111
$ws.Cells.Item(34, 10).Value = 2876.125
$ws.Cells.Item(34, 11).Value = 2841.111
$ws.Cells.Item(34, 12).Value = 2876.125
$ws.Cells.Item(34, 13).Value = -2639.111
$ws.Cells.Item(34, 14).Value = -3280.125
$ws.Cells.Item(99, 8).Value = 2849.3635
$ws.Cells.Item(99, 10).Value = 4124.75
$ws.Cells.Item(99, 12).Value = 4124.75
$ws.Cells.Item(99, 14).Value = -7120.75
$ws.Cells.Item(107, 8).Value = 775.13336
$ws.Cells.Item(107, 9).Value = 473.42856
$ws.Cells.Item(107, 10).Value = 4999
$ws.Cells.Item(107, 11).Value = 473.42856
$ws.Cells.Item(107, 12).Value = 4999
$ws.Cells.Item(107, 13).Value = 1446.57144
$ws.Cells.Item(107, 14).Value = -8839
$ws.Cells.Item(126, 8).Value = 2849.3635
$ws.Cells.Item(126, 10).Value = 4124.75
$ws.Cells.Item(126, 12).Value = 12374.25
$ws.Cells.Item(126, 14).Value = -17314.25

# CUL sheet (index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132, 8).Value = 1214.8572
$ws.Cells.Item(132, 10).Value = 1424.75
$ws.Cells.Item(132, 12).Value = 12822.75
$ws.Cells.Item(132, 14).Value = -17882.75

# GSM sheet (index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2283.8333
$ws.Cells.Item(102, 9).Value = 2132.1875
$ws.Cells.Item(102, 11).Value = 2132.1875
$ws.Cells.Item(102, 13).Value = -510.1875
$ws.Cells.Item(113, 8).Value = 1105.6154
$ws.Cells.Item(113, 9).Value = 870.4
$ws.Cells.Item(113, 10).Value = 1252.625
$ws.Cells.Item(113, 11).Value = 870.4
$ws.Cells.Item(113, 12).Value = 1252.625
$ws.Cells.Item(113, 13).Value = 1299.6
$ws.Cells.Item(113, 14).Value = -5592.625
$ws.Cells.Item(132, 8).Value = 2140791.5
$ws.Cells.Item(132, 9).Value = 2962127
$ws.Cells.Item(132, 11).Value = 8886381
$ws.Cells.Item(132, 13).Value = -8883851
$ws.Cells.Item(134, 8).Value = 28849.666
$ws.Cells.Item(134, 10).Value = 28849.666
$ws.Cells.Item(134, 12).Value = 86548.99800000001
$ws.Cells.Item(134, 14).Value = -91618.99800000001

# LTW sheet (index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 2083.6956
$ws.Cells.Item(7, 9).Value = 1809.5294
$ws.Cells.Item(7, 10).Value = 2860.5
$ws.Cells.Item(7, 11).Value = 1809.5294
$ws.Cells.Item(7, 12).Value = 2860.5
$ws.Cells.Item(7, 13).Value = -1697.5294
$ws.Cells.Item(7, 14).Value = -3084.5
$ws.Cells.Item(40, 8).Value = 11632.8
$ws.Cells.Item(40, 9).Value = 15387.625
$ws.Cells.Item(40, 10).Value = 7341.5713
$ws.Cells.Item(40, 11).Value = 15387.625
$ws.Cells.Item(40, 12).Value = 7341.5713
$ws.Cells.Item(40, 13).Value = -15251.625
$ws.Cells.Item(40, 14).Value = -7613.5713
$ws.Cells.Item(46, 8).Value = 1603.125
$ws.Cells.Item(46, 9).Value = 906.9
$ws.Cells.Item(46, 10).Value = 2763.5
$ws.Cells.Item(46, 11).Value = 906.9
$ws.Cells.Item(46, 12).Value = 2763.5
$ws.Cells.Item(46, 13).Value = -718.9
$ws.Cells.Item(46, 14).Value = -3139.5
$ws.Cells.Item(47, 8).Value = 500005000
$ws.Cells.Item(47, 10).Value = 9999
$ws.Cells.Item(47, 12).Value = 9999
$ws.Cells.Item(47, 14).Value = -10979
$ws.Cells.Item(48, 8).Value = 25055
$ws.Cells.Item(48, 9).Value = 25055
$ws.Cells.Item(48, 11).Value = 25055
$ws.Cells.Item(48, 13).Value = -24394
$ws.Cells.Item(52, 8).Value = 500005000
$ws.Cells.Item(52, 10).Value = 9999
$ws.Cells.Item(52, 12).Value = 9999
$ws.Cells.Item(52, 14).Value = -10465
$ws.Cells.Item(68, 8).Value = 2005.8889
$ws.Cells.Item(68, 10).Value = 4500
$ws.Cells.Item(68, 12).Value = 4500
$ws.Cells.Item(68, 14).Value = -5998
$ws.Cells.Item(71, 8).Value = 2005.8889
$ws.Cells.Item(71, 10).Value = 4500
$ws.Cells.Item(71, 12).Value = 22500
$ws.Cells.Item(71, 14).Value = -29988
$ws.Cells.Item(82, 8).Value = 1686.1666
$ws.Cells.Item(82, 9).Value = 1226.4
$ws.Cells.Item(82, 10).Value = 3985
$ws.Cells.Item(82, 11).Value = 1226.4
$ws.Cells.Item(82, 12).Value = 3985
$ws.Cells.Item(82, 13).Value = -865.4000000000001
$ws.Cells.Item(82, 14).Value = -4707
$ws.Cells.Item(85, 8).Value = 1686.1666
$ws.Cells.Item(85, 9).Value = 1226.4
$ws.Cells.Item(85, 10).Value = 3985
$ws.Cells.Item(85, 11).Value = 1226.4
$ws.Cells.Item(85, 12).Value = 3985
$ws.Cells.Item(85, 13).Value = 21.59999999999991
$ws.Cells.Item(85, 14).Value = -6481
$ws.Cells.Item(93, 8).Value = 14493604
$ws.Cells.Item(93, 10).Value = 47620016
$ws.Cells.Item(93, 12).Value = 47620016
$ws.Cells.Item(93, 14).Value = -47622512
$ws.Cells.Item(122, 8).Value = 12519.8
$ws.Cells.Item(122, 9).Value = 9533
$ws.Cells.Item(122, 11).Value = 28599
$ws.Cells.Item(122, 13).Value = -26149
$ws.Cells.Item(126, 8).Value = 2083.6956
$ws.Cells.Item(126, 9).Value = 1809.5294
$ws.Cells.Item(126, 10).Value = 2860.5
$ws.Cells.Item(126, 11).Value = 5428.5882
$ws.Cells.Item(126, 12).Value = 8581.5
$ws.Cells.Item(126, 13).Value = -2958.5882
$ws.Cells.Item(126, 14).Value = -13521.5
$ws.Cells.Item(132, 8).Value = 2107.4644
$ws.Cells.Item(132, 9).Value = 1432.5
$ws.Cells.Item(132, 10).Value = 2377.45
$ws.Cells.Item(132, 11).Value = 4297.5
$ws.Cells.Item(132, 12).Value = 7132.349999999999
$ws.Cells.Item(132, 13).Value = -1767.5
$ws.Cells.Item(132, 14).Value = -12192.35
$ws.Cells.Item(136, 8).Value = 3596.7932
$ws.Cells.Item(136, 9).Value = 2938.923
$ws.Cells.Item(136, 11).Value = 8816.769
$ws.Cells.Item(136, 13).Value = -6266.769

# WVR sheet (index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 1508.2
$ws.Cells.Item(132, 9).Value = 997.069
$ws.Cells.Item(132, 11).Value = 2991.207
$ws.Cells.Item(132, 13).Value = -461.2069999999999
